$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 424.23914
$ws.Range("I9").Value = 462.025
$ws.Range("J9").Value = 172.33333
$ws.Range("K9").Value = 462.025
$ws.Range("L9").Value = 172.33333
$ws.Range("M9").Value = -293.025
$ws.Range("N9").Value = -510.33333
$ws.Range("H13").Value = 2000
$ws.Range("I13").Value = 3000
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 3000
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = -2831
$ws.Range("N13").Value = -1338
$ws.Range("H43").Value = 483.66666
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 483.66666
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 483.66666
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -621.66666
$ws.Range("H69").Value = 11366664
$ws.Range("I69").Value = 41668668
$ws.Range("J69").Value = 3412.5
$ws.Range("K69").Value = 125006004
$ws.Range("L69").Value = 10237.5
$ws.Range("M69").Value = -125005130
$ws.Range("N69").Value = -11985.5
$ws.Range("H72").Value = 11366664
$ws.Range("I72").Value = 41668668
$ws.Range("J72").Value = 3412.5
$ws.Range("K72").Value = 375018012
$ws.Range("L72").Value = 30712.5
$ws.Range("M72").Value = -375013644
$ws.Range("N72").Value = -39448.5
$ws.Range("H129").Value = 946.5238000000001
$ws.Range("I129").Value = 306.7143
$ws.Range("J129").Value = 1266.4286
$ws.Range("K129").Value = 920.1428999999999
$ws.Range("L129").Value = 3799.2858
$ws.Range("M129").Value = 4079.8571
$ws.Range("N129").Value = -13799.2858
$ws.Range("H132").Value = 7584.1113
$ws.Range("I132").Value = 8367.200000000001
$ws.Range("J132").Value = 3668.6667
$ws.Range("K132").Value = 25101.6
$ws.Range("L132").Value = 11006.0001
$ws.Range("M132").Value = -22571.6
$ws.Range("N132").Value = -16066.0001
$ws.Range("H138").Value = 296093.4
$ws.Range("I138").Value = 455359.3
$ws.Range("J138").Value = 4105.9165
$ws.Range("K138").Value = 1366077.9
$ws.Range("L138").Value = 12317.7495
$ws.Range("M138").Value = -1360937.9
$ws.Range("N138").Value = -22597.7495
$ws.Range("H141").Value = 4957.971
$ws.Range("I141").Value = 1829.5927
$ws.Range("J141").Value = 15516.25
$ws.Range("K141").Value = 5488.7781
$ws.Range("L141").Value = 46548.75
$ws.Range("M141").Value = -308.7780999999995
$ws.Range("N141").Value = -56908.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 8000
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H92").Value = 26930
$ws.Range("J92").Value = 26930
$ws.Range("L92").Value = 26930
$ws.Range("N92").Value = -31922
$ws.Range("H110").Value = 7661.421
$ws.Range("I110").Value = 9170.571
$ws.Range("K110").Value = 9170.571
$ws.Range("M110").Value = -7125.571

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H88").Value = 34800
$ws.Range("J88").Value = 34800
$ws.Range("L88").Value = 34800
$ws.Range("N88").Value = -35612
$ws.Range("H91").Value = 34800
$ws.Range("J91").Value = 34800
$ws.Range("L91").Value = 34800
$ws.Range("N91").Value = -37608

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1453.8928
$ws.Range("I58").Value = 712.9375
$ws.Range("K58").Value = 712.9375
$ws.Range("M58").Value = -509.9375
$ws.Range("H62").Value = 90911900
$ws.Range("J62").Value = 166669330
$ws.Range("L62").Value = 166669330
$ws.Range("N62").Value = -166670578
$ws.Range("H65").Value = 90911900
$ws.Range("J65").Value = 166669330
$ws.Range("L65").Value = 833346650
$ws.Range("N65").Value = -833352890
$ws.Range("H136").Value = 1453.8928
$ws.Range("I136").Value = 712.9375
$ws.Range("K136").Value = 2138.8125
$ws.Range("M136").Value = 411.1875

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 60063.332
$ws.Range("J2").Value = 61.090908
$ws.Range("L2").Value = 366.545448
$ws.Range("N2").Value = -592.545448
$ws.Range("H17").Value = 280.2
$ws.Range("I17").Value = 225
$ws.Range("J17").Value = 501
$ws.Range("K17").Value = 675
$ws.Range("L17").Value = 1503
$ws.Range("M17").Value = -506
$ws.Range("N17").Value = -1841
$ws.Range("H34").Value = 298.24243
$ws.Range("I34").Value = 120.8
$ws.Range("J34").Value = 446.1111
$ws.Range("K34").Value = 362.4
$ws.Range("L34").Value = 1338.3333
$ws.Range("M34").Value = -278.4
$ws.Range("N34").Value = -1506.3333
$ws.Range("H39").Value = 2312.5
$ws.Range("I39").Value = 1000
$ws.Range("J39").Value = 2500
$ws.Range("K39").Value = 3000
$ws.Range("L39").Value = 7500
$ws.Range("M39").Value = -2706
$ws.Range("N39").Value = -8088
$ws.Range("H64").Value = 40004456
$ws.Range("I64").Value = 2500
$ws.Range("J64").Value = 50004944
$ws.Range("K64").Value = 7500
$ws.Range("L64").Value = 150014832
$ws.Range("M64").Value = -7230
$ws.Range("N64").Value = -150015372
$ws.Range("H67").Value = 40004456
$ws.Range("I67").Value = 2500
$ws.Range("J67").Value = 50004944
$ws.Range("K67").Value = 7500
$ws.Range("L67").Value = 150014832
$ws.Range("M67").Value = -6564
$ws.Range("N67").Value = -150016704

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 3858
$ws.Range("I9").Value = 787
$ws.Range("J9").Value = 10000
$ws.Range("K9").Value = 787
$ws.Range("L9").Value = 10000
$ws.Range("M9").Value = -617
$ws.Range("N9").Value = -10340
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H125").Value = 53280
$ws.Range("J125").Value = 53280
$ws.Range("L125").Value = 53280
$ws.Range("N125").Value = -58200
$ws.Range("H132").Value = 2055.5173
$ws.Range("I132").Value = 1661.3334
$ws.Range("J132").Value = 3090.25
$ws.Range("K132").Value = 4984.0002
$ws.Range("L132").Value = 9270.75
$ws.Range("M132").Value = -2454.0002
$ws.Range("N132").Value = -14330.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 18101
$ws.Range("I12").Value = 500
$ws.Range("J12").Value = 23968
$ws.Range("K12").Value = 500
$ws.Range("L12").Value = 23968
$ws.Range("M12").Value = -330
$ws.Range("N12").Value = -24308
$ws.Range("H46").Value = 797.5714
$ws.Range("J46").Value = 797.6667
$ws.Range("L46").Value = 797.6667
$ws.Range("N46").Value = -1173.6667
$ws.Range("H55").Value = 291.64285
$ws.Range("I55").Value = 310.375
$ws.Range("J55").Value = 266.66666
$ws.Range("K55").Value = 310.375
$ws.Range("L55").Value = 266.66666
$ws.Range("M55").Value = -137.375
$ws.Range("N55").Value = -612.66666
$ws.Range("H136").Value = 3792911
$ws.Range("I136").Value = 6952327.5
$ws.Range("J136").Value = 1611.3334
$ws.Range("K136").Value = 20856982.5
$ws.Range("L136").Value = 4834.0002
$ws.Range("M136").Value = -20854432.5
$ws.Range("N136").Value = -9934.0002

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 100002770
$ws.Range("I81").Value = 166669200
$ws.Range("J81").Value = 3125
$ws.Range("K81").Value = 333338400
$ws.Range("L81").Value = 6250
$ws.Range("M81").Value = -333337339
$ws.Range("N81").Value = -8372
$ws.Range("H84").Value = 100002770
$ws.Range("I84").Value = 166669200
$ws.Range("J84").Value = 3125
$ws.Range("K84").Value = 1666692000
$ws.Range("L84").Value = 31250
$ws.Range("M84").Value = -1666686696
$ws.Range("N84").Value = -41858
